$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add a new "instrument block" (rows 21-27) that mirrors the existing
# 交易品种/交易数量/止盈目标/级别/1D/30F/5F block found at rows 12-18,
# for the new instrument 科创50（588000） with targets 1.452/1.464.
# ---------------------------------------------------------------------------

# Row 21: 交易品种： / 科创50（588000）  -- copy shape+style from row 12, then
# replace the instrument name.
$ws.Range("A12:B12").Copy($ws.Range("A21:B21"))
$ws.Range("B21").Value = "科创50（588000）"

# Row 22: 交易数量： / 无  -- identical to row 13.
$ws.Range("A13:B13").Copy($ws.Range("A22:B22"))

# Row 23: 止盈目标： / 翻倍  -- identical to row 14.
$ws.Range("A14:B14").Copy($ws.Range("A23:B23"))

# Row 24: header row (级别 / 中枢0底/中枢0顶 / 中枢0执行情况 / 中枢1底/中枢1顶 / 中枢1执行情况)
# -- identical to row 15.
$ws.Range("A15:R15").Copy($ws.Range("A24:R24"))

# Row 25: 1D row -- same styling as row 16, but no value recorded yet.
$ws.Range("A16:R16").Copy($ws.Range("A25:R25"))
$ws.Range("B25").ClearContents()

# Row 26: 30F row -- same styling as row 17, but no value recorded yet.
$ws.Range("A17:R17").Copy($ws.Range("A26:R26"))
$ws.Range("B26").ClearContents()

# Row 27: 5F row -- same styling as row 18, with the new target values.
# (Row 18 also carries a full-row fill beyond column H; AutoFit below just
# clears the stray leftover row height that the placeholder row 27 had.)
$ws.Range("A18:H18").Copy($ws.Range("A27:H27"))
$ws.Range("B27").Value = "1.452/1.464"
$ws.Range("D27").ClearContents()
$ws.Range("E27").ClearContents()
$ws.Rows.Item(27).AutoFit()

# ---------------------------------------------------------------------------
# Update the active selection to reflect where the author left off editing.
# ---------------------------------------------------------------------------
$ws.Range("B31").Select()
